$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B (id) and C (speaker_variant) values for rows 2-29 to reflect the
# re-exported grouping (no preference / no levenshtein distance logic).
$ws.Range("B2").Value = "#ka"
$ws.Range("C2").Value = "Ka"
$ws.Range("B3").Value = "#alcinea"
$ws.Range("C3").Value = "Alcinea"
$ws.Range("B4").Value = "#ni"
$ws.Range("C4").Value = "Ni"
$ws.Range("B5").Value = "#kar"
$ws.Range("C5").Value = "Kar"
$ws.Range("B6").Value = "#poli"
$ws.Range("C6").Value = "Poli"
$ws.Range("B7").Value = "#pol"
$ws.Range("C7").Value = "Pol"
$ws.Range("B8").Value = "#al"
$ws.Range("C8").Value = "Al"
$ws.Range("B9").Value = "#kla"
$ws.Range("C9").Value = "Kla"
$ws.Range("B10").Value = "#sid"
$ws.Range("C10").Value = "Sid"
$ws.Range("B11").Value = "#klarimeen"
$ws.Range("C11").Value = "Klarimeen"
$ws.Range("B12").Value = "#ha"
$ws.Range("C12").Value = "Ha"
$ws.Range("B13").Value = "#kl"
$ws.Range("C13").Value = "Kl"
$ws.Range("B14").Value = "#rey"
$ws.Range("C14").Value = "Rey"
$ws.Range("B15").Value = "#tol"
$ws.Range("C15").Value = "Tol"
$ws.Range("B16").Value = "#cla"
$ws.Range("C16").Value = "Cla"
$ws.Range("B17").Value = "#polin"
$ws.Range("C17").Value = "Polin"
$ws.Range("B18").Value = "#si"
$ws.Range("C18").Value = "Si"
$ws.Range("B19").Value = "#ag"
$ws.Range("C19").Value = "Ag"
$ws.Range("B20").Value = "#alb"
$ws.Range("C20").Value = "Alb"
$ws.Range("B21").Value = "#art"
$ws.Range("C21").Value = "Art"
$ws.Range("B22").Value = "#phi"
$ws.Range("C22").Value = "Phi"
$ws.Range("B23").Value = "#mar"
$ws.Range("C23").Value = "Mar"
$ws.Range("B24").Value = "#kl,"
$ws.Range("C24").Value = "Kl,"
$ws.Range("B25").Value = "#oct"
$ws.Range("C25").Value = "Oct"
$ws.Range("B26").Value = "#po"
$ws.Range("C26").Value = "Po"
$ws.Range("B27").Value = "#pa"
$ws.Range("C27").Value = "Pa"
$ws.Range("B28").Value = "#nis"
$ws.Range("C28").Value = "Nis"
$ws.Range("B29").Value = "#cl"
$ws.Range("C29").Value = "Cl"

# Clear the "is_prefered" (D) column for rows 2-10, which previously held "x".
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("D10").Value = ""
